# Insert a new weekly data row at row 333 (pushing existing rows 333-359
# down to 334-360), then populate the new row with the latest observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(333).Insert()

$ws.Cells.Item(333, 1).Value  = 6
$ws.Cells.Item(333, 2).Value  = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(333, 3).Value  = 'Metropolitana'
$ws.Cells.Item(333, 4).Value  = 44578
$ws.Cells.Item(333, 5).Value  = 13
$ws.Cells.Item(333, 6).Value  = 100112039
$ws.Cells.Item(333, 7).Value  = 'Ciboulette'
$ws.Cells.Item(333, 8).Value  = 'Sin especificar'
$ws.Cells.Item(333, 9).Value  = 'Primera'
$ws.Cells.Item(333, 10).Value = 580
$ws.Cells.Item(333, 11).Value = 1000
$ws.Cells.Item(333, 12).Value = 1200
$ws.Cells.Item(333, 13).Value = 1086
$ws.Cells.Item(333, 14).Value = '$/docena de atados'
$ws.Cells.Item(333, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(333, 16).Value = 362
$ws.Cells.Item(333, 17).Value = 3
$ws.Cells.Item(333, 18).Value = 'Hortaliza'
